$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 41445
$ws.Range("D2").Value = 59849291
$ws.Range("C3").Value = 98443
$ws.Range("D3").Value = 144211142
$ws.Range("C4").Value = 33494
$ws.Range("D4").Value = 49577994
$ws.Range("C5").Value = 9562
$ws.Range("D5").Value = 14203242
$ws.Range("C6").Value = 2357
$ws.Range("D6").Value = 3500973
$ws.Range("C7").Value = 237
$ws.Range("D7").Value = 350593
$ws.Range("C12").Value = 44590
$ws.Range("D12").Value = 60352534
$ws.Range("C13").Value = 10494
$ws.Range("D13").Value = 15156296
$ws.Range("C14").Value = 27721
$ws.Range("D14").Value = 40621145
$ws.Range("C15").Value = 8790
$ws.Range("D15").Value = 13043922
$ws.Range("C16").Value = 2325
$ws.Range("D16").Value = 3455103
$ws.Range("C17").Value = 481
$ws.Range("D17").Value = 710623
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 70500
$ws.Range("C20").Value = 10947
$ws.Range("D20").Value = 14418427
$ws.Range("C21").Value = 14428
$ws.Range("D21").Value = 20799935
$ws.Range("C22").Value = 33597
$ws.Range("D22").Value = 49258005
$ws.Range("C23").Value = 10809
$ws.Range("D23").Value = 16060937
$ws.Range("C24").Value = 2870
$ws.Range("D24").Value = 4265115
$ws.Range("C25").Value = 592
$ws.Range("D25").Value = 882092
$ws.Range("C27").Value = 12465
$ws.Range("D27").Value = 16572479
$ws.Range("C28").Value = 8403
$ws.Range("D28").Value = 12151828
$ws.Range("C29").Value = 24144
$ws.Range("D29").Value = 35418075
$ws.Range("C30").Value = 8310
$ws.Range("D30").Value = 12353707
$ws.Range("C31").Value = 2118
$ws.Range("D31").Value = 3159208
$ws.Range("C32").Value = 421
$ws.Range("D32").Value = 622121
$ws.Range("C33").Value = 37
$ws.Range("D33").Value = 55393
$ws.Range("C34").Value = 8955
$ws.Range("D34").Value = 11800564
$ws.Range("C35").Value = 3660
$ws.Range("D35").Value = 5282462
$ws.Range("C36").Value = 8596
$ws.Range("D36").Value = 12560045
$ws.Range("C37").Value = 3398
$ws.Range("D37").Value = 5038008
$ws.Range("C38").Value = 868
$ws.Range("D38").Value = 1293055
$ws.Range("C39").Value = 179
$ws.Range("D39").Value = 266186
$ws.Range("C41").Value = 2741
$ws.Range("D41").Value = 3692222
$ws.Range("C42").Value = 18757
$ws.Range("D42").Value = 27084869
$ws.Range("C43").Value = 54658
$ws.Range("D43").Value = 80076687
$ws.Range("C44").Value = 20056
$ws.Range("D44").Value = 29770953
$ws.Range("C45").Value = 6017
$ws.Range("D45").Value = 8948637
$ws.Range("C46").Value = 1402
$ws.Range("D46").Value = 2092644
$ws.Range("C50").Value = 18200
$ws.Range("D50").Value = 24102640
$ws.Range("C51").Value = 2356
$ws.Range("D51").Value = 3417726
$ws.Range("C52").Value = 7945
$ws.Range("D52").Value = 11666298
$ws.Range("C53").Value = 2662
$ws.Range("D53").Value = 3973133
$ws.Range("C54").Value = 837
$ws.Range("D54").Value = 1250414
$ws.Range("C57").Value = 7926
$ws.Range("D57").Value = 10906865
$ws.Range("C58").Value = 1631
$ws.Range("D58").Value = 3261618
$ws.Range("C59").Value = 3900
$ws.Range("D59").Value = 7768322
$ws.Range("C60").Value = 1533
$ws.Range("D60").Value = 3055462
$ws.Range("C61").Value = 509
$ws.Range("D61").Value = 1008083
$ws.Range("C64").Value = 2538
$ws.Range("D64").Value = 4696528
$ws.Range("C65").Value = 16919
$ws.Range("D65").Value = 24414710
$ws.Range("C66").Value = 48144
$ws.Range("D66").Value = 70375958
$ws.Range("C67").Value = 16824
$ws.Range("D67").Value = 24995976
$ws.Range("C68").Value = 4909
$ws.Range("D68").Value = 7311024
$ws.Range("C69").Value = 1076
$ws.Range("D69").Value = 1599699
$ws.Range("C73").Value = 16099
$ws.Range("D73").Value = 21130225
$ws.Range("C74").Value = 62123
$ws.Range("D74").Value = 90299949
$ws.Range("C75").Value = 169277
$ws.Range("D75").Value = 249128604
$ws.Range("C76").Value = 72223
$ws.Range("D76").Value = 107563242
$ws.Range("C77").Value = 23567
$ws.Range("D77").Value = 35196060
$ws.Range("C78").Value = 6079
$ws.Range("D78").Value = 9075795
$ws.Range("C79").Value = 442
$ws.Range("D79").Value = 657460
$ws.Range("C85").Value = 61101
$ws.Range("D85").Value = 82450817
$ws.Range("C86").Value = 5126
$ws.Range("D86").Value = 7425373
$ws.Range("C87").Value = 12587
$ws.Range("D87").Value = 18483665
$ws.Range("C88").Value = 4128
$ws.Range("D88").Value = 6151126
$ws.Range("C89").Value = 1451
$ws.Range("D89").Value = 2167611
$ws.Range("C90").Value = 348
$ws.Range("D90").Value = 518512
$ws.Range("C93").Value = 5857
$ws.Range("D93").Value = 7852374
$ws.Range("C94").Value = 1829
$ws.Range("D94").Value = 2632297
$ws.Range("C95").Value = 5825
$ws.Range("D95").Value = 8583552
$ws.Range("C96").Value = 2099
$ws.Range("D96").Value = 3128431
$ws.Range("C98").Value = 218
$ws.Range("D98").Value = 329113
$ws.Range("C99").Value = 23
$ws.Range("D99").Value = 34500
$ws.Range("C101").Value = 3966
$ws.Range("D101").Value = 5257873
$ws.Range("C102").Value = 912
$ws.Range("D102").Value = 1758061
$ws.Range("C103").Value = 610
$ws.Range("D103").Value = 1244284
$ws.Range("C107").Value = 11864
$ws.Range("D107").Value = 17195998
$ws.Range("C108").Value = 31068
$ws.Range("D108").Value = 45599236
$ws.Range("C109").Value = 10401
$ws.Range("D109").Value = 15461927
$ws.Range("C110").Value = 2883
$ws.Range("D110").Value = 4297571
$ws.Range("C111").Value = 563
$ws.Range("D111").Value = 838453
$ws.Range("C112").Value = 66
$ws.Range("D112").Value = 99000
$ws.Range("C115").Value = 10466
$ws.Range("D115").Value = 13774790
$ws.Range("C116").Value = 33292
$ws.Range("D116").Value = 47962665
$ws.Range("C117").Value = 70729
$ws.Range("D117").Value = 103450163
$ws.Range("C118").Value = 22700
$ws.Range("D118").Value = 33718034
$ws.Range("C119").Value = 6507
$ws.Range("D119").Value = 9687732
$ws.Range("C120").Value = 1293
$ws.Range("D120").Value = 1931237
$ws.Range("C121").Value = 119
$ws.Range("D121").Value = 174895
$ws.Range("C125").Value = 27623
$ws.Range("D125").Value = 36780440
$ws.Range("C126").Value = 39778
$ws.Range("D126").Value = 57344295
$ws.Range("C127").Value = 83066
$ws.Range("D127").Value = 121368113
$ws.Range("C128").Value = 25454
$ws.Range("D128").Value = 37770544
$ws.Range("C129").Value = 6943
$ws.Range("D129").Value = 10317009
$ws.Range("C130").Value = 1466
$ws.Range("D130").Value = 2172096
$ws.Range("C134").Value = 34132
$ws.Range("D134").Value = 45195433
$ws.Range("C135").Value = 14477
$ws.Range("D135").Value = 20943320
$ws.Range("C136").Value = 34511
$ws.Range("D136").Value = 50654972
$ws.Range("C137").Value = 12174
$ws.Range("D137").Value = 18087573
$ws.Range("C138").Value = 3227
$ws.Range("D138").Value = 4810375
$ws.Range("C139").Value = 583
$ws.Range("D139").Value = 868490
$ws.Range("C143").Value = 11567
$ws.Range("D143").Value = 15371757
$ws.Range("C144").Value = 39020
$ws.Range("D144").Value = 56341613
$ws.Range("C145").Value = 89292
$ws.Range("D145").Value = 130734084
$ws.Range("C146").Value = 26590
$ws.Range("D146").Value = 39486431
$ws.Range("C147").Value = 7076
$ws.Range("D147").Value = 10544735
$ws.Range("C148").Value = 1684
$ws.Range("D148").Value = 2499966
$ws.Range("C149").Value = 116
$ws.Range("D149").Value = 173630
$ws.Range("C150").Value = 19
$ws.Range("D150").Value = 28500
$ws.Range("C151").Value = 31605
$ws.Range("D151").Value = 42484434
